# Daily attendance processing - 2026-01-03 19:30:03
# Normalizes the "Recorded By" (column G) entries: for a fixed set of
# known exact text values, rotate the comma-separated list of names so
# that the last entry moves to the front (e.g. "System, X" -> "X, System").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

# Only these exact "Recorded By" strings are affected by this pass.
$targets = @(
    "System, dnasr281@gmail.com",
    "System, admin@admin.com",
    "System, backup@backdoor.com, system"
)

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($targets -contains $val) {
        $parts = $val -split ", "
        $n = $parts.Count
        $rotated = @($parts[$n - 1]) + $parts[0..($n - 2)]
        $newVal = $rotated -join ", "
        $cell.Value2 = $newVal
    }
}
